# Bump the "Förändrad" (last-changed) date in column C from 2023-09-20
# (serial 45189) to 2023-09-21 (serial 45190) for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45189) {
        $cell.Value2 = 45190
    }
}
